$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row at position 14 (Brane's "Spikefield Hazard" pickup).
#    This pushes the existing rows 14-27 down to 15-28, and Excel adjusts the
#    existing SUM()/shared formulas in column F/G automatically.
# ---------------------------------------------------------------------------
$ws.Rows(14).Insert()

$ws.Range("A14").Value = "Spikefield Hazard: Spikefield Cave"
$ws.Range("B14").Value = "Zendikar Rising"
$ws.Range("C14").Value = "Normal"
$ws.Range("D14").Value = 0.37
$ws.Range("E14").Value = 1
$ws.Range("F14").Formula = "=D14*E14"

# ---------------------------------------------------------------------------
# 2) Updated prices (column D) for existing cards - rows numbered as they are
#    *after* the insert above shifted everything from row 14 onward down by 1.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 2.97     # Arclight Phoenix
$ws.Range("D3").Value = 0.28     # Chart a Course
$ws.Range("D4").Value = 0.82     # Consider
$ws.Range("D6").Value = 4.46     # Hall of Storm Giants
$ws.Range("D7").Value = 0.3      # Fiery Impulse
$ws.Range("D8").Value = 0.18     # Izzet Charm
$ws.Range("D9").Value = 0.12     # Lightning Axe
$ws.Range("D10").Value = 0.09    # Opt
$ws.Range("D12").Value = 3.14    # Search for Azcanta: Azcanta the Sunken Ruin
$ws.Range("D13").Value = 1.06    # Shivan Reef
$ws.Range("D15").Value = 13.84   # Steam Vents
$ws.Range("D16").Value = 1.63    # Sulfur Falls
$ws.Range("D17").Value = 0.27    # Temple of Epiphany
$ws.Range("D18").Value = 3.6     # Thing in the Ice: Awoken Horror
$ws.Range("D19").Value = 0.16    # Treasure Cruise
$ws.Range("D20").Value = 0.06    # Flame-Blessed Bolt
$ws.Range("D21").Value = 0.28    # Invasive Surgery
$ws.Range("D22").Value = 0.18    # Lava Coil
$ws.Range("D23").Value = 0.5     # Mystical Dispute
$ws.Range("D24").Value = 0.76    # Narset, Parter of Veils
$ws.Range("D26").Value = 0.03    # Negate (Magic 2014 copy)
$ws.Range("D27").Value = 0.2     # Sweltering Suns
